# commands: sync process from package.json to commands
#
# Reorders / rewrites the "developer" command list (rows 12-25) on the
# "application deploy commands" sheet so that it matches the steps actually
# used in package.json, and tidies up the trailing rows that are no longer
# needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: now the "cd <source>" step (used to be the npm build step) ---
$ws.Range("D12").Formula = '="cd "&$B$5'

# --- Row 13: now the "cp -r backend <version>" step (used to be npm sync) ---
$ws.Range("D13").Formula = '=" ""cp -r backend "" & $B$4"'
$ws.Range("D13").Formula = '= "cp -r backend " & $B$4'

# --- Row 14: now the npm build step (used to be the "cd <source>" step) ---
$ws.Range("D14").Formula = '="npm run build --prefix "&$B$5 & "/frontend"'

# --- Row 15: brand-new multi-command rsync/touch formula, replacing the old
#     "cp -r backend" formula that used to live here ---
$ws.Range("D15").Formula = '="touch "&$B$5&"/" & $B$4& "/base/static/empty.txt && rm -r "&$B$5&"/" & $B$4& "/base/static/* && rsync -av " & $B$5 &"/frontend/build/index.html "&$B$5&"/" & $B$4& "/base/templates/. && rsync -av --exclude=''static'' "  & $B$5 &"/frontend/build/* "& $B$5 &"/" & $B$4& "/base/static/. && rsync -av " & $B$5 &"/frontend/build/static/* "&$B$5&"/" & $B$4& "/base/static/."'

# Rows 16 (source venv) and 17 (pip freeze) are unchanged.

# --- Row 18: now the "python manage.py collectstatic" step (used to be "cd
#     <source>/<version>") ---
$ws.Range("D18").Formula = '="python "&$B$5&"/"&$B$4&"/manage.py collectstatic"'

# --- Row 19: now just the literal "deactivate" shared string (used to be the
#     collectstatic formula) ---
$ws.Range("D19").Value = "deactivate"

# --- Row 20: now the "vi .../backend/settings.py" step (used to be the bare
#     "deactivate" text) ---
$ws.Range("D20").Formula = '="vi "&$B$5&"/"&$B$4&"/backend/settings.py"'

# --- Row 21: used to hold the "vi settings.py" step; now blank again ---
$ws.Range("C21:D21").ClearContents()

# --- Row 22: used to be a blank spacer row; now holds the "cd <source>" step ---
$ws.Range("C22").Value = "developer"
$ws.Range("D22").Formula = '="cd "&$B$5'

# --- Row 23: now the "zip -r" step (used to be "cd <source>") ---
$ws.Range("D23").Formula = '="zip -r "&$B$4&".zip "&$B$4'

# --- Row 24: now the "rm -r <version>" step (used to be the "zip -r" step) ---
$ws.Range("D24").Formula = '="rm -r "&B4'

# --- Row 25: used to hold the "rm -r <version>" step; now removed entirely ---
$ws.Range("C25:D25").ClearContents()

# Update the sheet's saved selection/scroll position to match.
$ws.Range("C22").Select()
